$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# Sheet 1: "Indicadores"
# -----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Indicadores")

# Row labels (column A) and data (columns B..F) for rows 2..14
$sheet1Data = @(
    @{ Row = 2;  Label = "Endividamento Total";                     Values = @(2.61068036787363, 2.65978964066897, 2.590901228701166, 2.547669704768785, 2.535321483732147) },
    @{ Row = 3;  Label = "Dependência Financeira (%)";               Values = @(72.30438869921596, 72.67602517675823, 72.15183776130478, 71.81248303200837, 71.71402927282512) },
    @{ Row = 4;  Label = "Independência Financeira (%)";             Values = @(27.69561130078404, 27.32397482324178, 27.84816223869519, 28.18751696799164, 28.28597072717488) },
    @{ Row = 5;  Label = "Capital Terceiros LP";                     Values = @(1.777338748265412, 1.640271672050957, 1.451796019768496, 1.341621070215703, 1.214905401540097) },
    @{ Row = 6;  Label = "Imobilização Capital Permanente";          Values = @(0.5436158758912167, 0.5841007837308343, 0.6332860718162319, 0.6642000495078371, 0.7170454445870662) },
    @{ Row = 7;  Label = "Imobilização PL";                          Values = @(0.6924936979743261, 0.6298966221302366, 0.6069636238741689, 0.6115588787180845, 0.6117320399472587) },
    @{ Row = 8;  Label = "Imobilização Recursos Não Correntes";      Values = @(0.506399173702311, 0.5413007481956502, 0.5737451953735432, 0.592655594514736, 0.619949725220023) },
    @{ Row = 9;  Label = "Giro do Imobilizado";                      Values = @(1.707850277633619, 2.214638267166006, 2.517534106021249, 2.912558117819684, 2.796014152384233) },
    @{ Row = 10; Label = "Nível Automação (Imobilizado/Receita)";    Values = @(0.5855314210479801, 0.4515410100267366, 0.3972140824659633, 0.343340788251323, 0.3576519808196515) },
    @{ Row = 11; Label = "Liquidez Geral";                           Values = @(0.8443156610764809, 0.8386409868815722, 0.8430252061048135, 0.8477923465624001, 0.8528273445922694) },
    @{ Row = 12; Label = "Composição Endividamento LP (%)";          Values = @(68.07952325902824, 61.6692255271138, 56.03440238037519, 52.6607145229394, 47.91918537098822) },
    @{ Row = 13; Label = "Participação ANC no Ativo (%)";            Values = @(41.81498450315775, 42.13861736122597, 43.2395148561486, 43.84018130761216, 44.92343442240518) },
    @{ Row = 14; Label = "Alavancagem Financeira (PNC/PL)";          Values = @(1.777338748265412, 1.640271672050957, 1.451796019768496, 1.341621070215703, 1.214905401540097) }
)

# Rows 12-14 are brand new (beyond the original A1:F11 extent) and, unlike
# the pre-existing rows 2-11, don't automatically carry the bold/centered/
# bordered label style (s="1"). Clone formatting from an existing labeled
# row (row 2) onto those new label cells BEFORE writing their text, so the
# label text we write next is not clobbered by the format copy.
$labelFormatSource1 = $ws1.Cells.Item(2, 1)
foreach ($r in 12..14) {
    $labelFormatSource1.Copy($ws1.Cells.Item($r, 1))
}

foreach ($entry in $sheet1Data) {
    $r = $entry.Row
    $ws1.Cells.Item($r, 1).Value = $entry.Label
    for ($i = 0; $i -lt 5; $i++) {
        $ws1.Cells.Item($r, 2 + $i).Value = $entry.Values[$i]
    }
}

# -----------------------------------------------------------------------
# Sheet 2: "Dados Base"
# -----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Dados Base")

# The final layout inserts three new rows after "AtivoCirculante" (for the
# non-current-asset breakdown) and one new row before "LucroBruto" (for
# "DepreciacaoAmortizacao"), pushing every following row down. Rather than
# using Rows.Insert() (which shifts cells but also leaves stray/unused
# auto-derived styles behind), just rewrite every row 2..17 directly with
# its final label + values, since all of sheet 2's data rows change
# position and/or content anyway.

# Full data for rows 2..17 in the final layout.
$sheet2Data = @(
    @{ Row = 2;  Label = "AtivoCirculante";          Values = @(40549746.68045199, 37967991.53685346, 32876200.42487912, 30670519.513, 40326906) },
    @{ Row = 3;  Label = "AtivoNaoCirculante";        Values = @(29141300.63507662, 27650888.97605425, 25044728.79473388, 23942405.2206, 32892812) },
    @{ Row = 4;  Label = "AtivoImobilizado";          Values = @(13366071.25367216, 11293848.89331861, 9790271.260123599, 9414354.1516, 12669526) },
    @{ Row = 5;  Label = "AtivoRealizavelLP";         Values = @(1995054.146192828, 2026121.285431994, 2354678.50797864, 2578961.465, 4454043) },
    @{ Row = 6;  Label = "AtivoTotal";                Values = @(69691047.31552862, 65618880.5129077, 57920929.21961301, 54612924.7336, 73219718) },
    @{ Row = 7;  Label = "PassivoCirculante";         Values = @(16084627.91634959, 18279637.44694333, 18373669.44346336, 18565945.7574, 27347016) },
    @{ Row = 8;  Label = "PassivoNaoCirculante";      Values = @(34305057.82322473, 29409556.67532437, 23417345.436912, 20652951.5502, 25161794) },
    @{ Row = 9;  Label = "PatrimonioLiquido";         Values = @(19301361.5759543, 17929686.39064, 16129914.33923764, 15394027.426, 20710908) },
    @{ Row = 10; Label = "ReceitaLiquida";            Values = @(22827248.50145473, 25011789.94273384, 24647341.8045608, 27419853.60827202, 35424174) },
    @{ Row = 11; Label = "DepreciacaoAmortizacao";    Values = @(-4751012.752529023, -2741056.434950025, -5534712.458766987, -3130259.287387266, -2610113) },
    @{ Row = 12; Label = "LucroBruto";                Values = @(2797304.76988113, 3904740.967449251, 4950725.889129187, 4729220.957990627, 6382234) },
    @{ Row = 13; Label = "LucroOperacional";          Values = @(-1953707.982647893, 1163684.532499225, -583986.5696378001, 1598961.670603361, 3772121) },
    @{ Row = 14; Label = "LucroLiquido";              Values = @(-4179466.994926325, -295266.5330622761, -1100740.706332947, 823687.2244588722, 1923831) },
    @{ Row = 15; Label = "PassivoTotal";              Values = @(50389685.73957432, 47689194.1222677, 41791014.88037536, 39218897.3076, 52508810) },
    @{ Row = 16; Label = "CapitalPermanente";         Values = @(53606419.39917903, 47339243.06596437, 39547259.77614964, 36046978.9762, 45872702) },
    @{ Row = 17; Label = "AtivoPermanente";           Values = @(13366071.25367216, 11293848.89331861, 9790271.260123599, 9414354.1516, 12669526) }
)

# Rows 14-17 are brand new (beyond the original A1:F13 extent) and don't
# automatically carry the bold/centered/bordered label style (s="1") used
# by the rest of column A. Re-clone that canonical label formatting from a
# known-good cell onto every column-A label cell BEFORE writing text/
# values (harmless no-op for rows that already have it), so all of them
# land on the same style index, matching the original workbook's
# consistent formatting (and so the format copy can't clobber our text).
$labelFormatSource2 = $ws2.Cells.Item(2, 1)
foreach ($entry in $sheet2Data) {
    $labelFormatSource2.Copy($ws2.Cells.Item($entry.Row, 1))
}

foreach ($entry in $sheet2Data) {
    $r = $entry.Row
    $ws2.Cells.Item($r, 1).Value = $entry.Label
    for ($i = 0; $i -lt 5; $i++) {
        $ws2.Cells.Item($r, 2 + $i).Value = $entry.Values[$i]
    }
}
